$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Bought" (D) column updates
$ws.Range("D4").Value = 10
$ws.Range("D5").Value = 10
$ws.Range("D19").Value = 10
$ws.Range("D22").Value = 10
$ws.Range("D27").Value = 3

# "Have" (C) column updates
$ws.Range("C16").Value = 4
$ws.Range("C17").Value = 4
$ws.Range("C18").Value = 8
$ws.Range("C20").Value = 2
$ws.Range("C21").Value = 4
$ws.Range("C23").Value = 2
$ws.Range("C28").Value = 3

# Move the active selection to match the author's final cursor position
$ws.Range("D28").Select() | Out-Null
